# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de sheets to reflect the newly generated
# handback report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-13 16:22:50"
$wsZh.Range("E4").Value = "2016-03-13 16:22:50"
$wsZh.Range("H3").Value = "2016-03-13 16:23:23"
$wsZh.Range("H4").Value = "2016-03-13 16:23:23"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-13 16:22:54"
$wsDe.Range("E4").Value = "2016-03-13 16:22:54"
$wsDe.Range("H3").Value = "2016-03-13 16:23:29"
$wsDe.Range("H4").Value = "2016-03-13 16:23:29"
